$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at A, shifting every existing column right by one.
$ws.Columns(1).Insert()

# Header label for the new "Match ID" column (row 3 is the visible header row).
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Match ID value (23) for every player data row (4-19) plus the summary row (20).
$ws.Range("A4:A19").Value = 23
$ws.Range("A4:A19").Font.Bold = $true

# Row 20 is hidden; temporarily reveal it so the write doesn't bake in a
# synthesized row height, then restore the hidden state.
$ws.Rows(20).Hidden = $false
$ws.Range("A20").Value = 23
$ws.Rows(20).Hidden = $true

# Restore the selection to the new Match ID column's data range.
$ws.Range("A3:A19").Select()
